# Generate Report for Handoff
#
# This script swaps the "fb153bbd" and "ca03050d" rows (row 6 / row 7) on the
# Overview, zh-cn and de-de sheets, updates fb153bbd's handoff status to
# "Ready for handoff" together with a fresh "Latest Handoff" timestamp, and
# fixes up the hyperlink display text so it still matches the row contents.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (columns: A FileName, B PathAndName, C Extension,
#                      D PublishURL, E zh-cn, F de-de, G Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 6 becomes the former row 7 (ca03050d) contents.
# (Column D is blank in both the before and after state, so it is left alone
# - assigning an empty string would delete the cell instead of keeping it.)
$ov.Range("A6").Value = "ca03050d-a2be-45e6-8869-2411d5b68e55.md"
$ov.Range("B6").Value = "e2e\ca03050d-a2be-45e6-8869-2411d5b68e55.md"
$ov.Range("C6").Value = ".md"
$ov.Range("E6").Value = "Ready for handoff"
$ov.Range("F6").Value = "Ready for handoff"
$ov.Range("G6").Value = "2016-09-06 04:03:59"

# Row 7 becomes the former row 6 (fb153bbd) contents, now marked handed off.
$ov.Range("A7").Value = "fb153bbd-eda5-4ccd-8490-bd45369ad1db.md"
$ov.Range("B7").Value = "e2e\fb153bbd-eda5-4ccd-8490-bd45369ad1db.md"
$ov.Range("C7").Value = ".md"
$ov.Range("E7").Value = "Ready for handoff"
$ov.Range("F7").Value = "Ready for handoff"
$ov.Range("G7").Value = "2016-09-06 04:11:11"

# Rebuild the hyperlinks so the display text matches the swapped rows while
# the link targets stay exactly where they were (rId6 keeps pointing at the
# fb153bbd github blob, rId7 keeps pointing at the ca03050d github blob).
$ovLinks = $ov.Hyperlinks
$ovTargets = @(
  @{ Cell = "B2"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5f8b17538436e140812b57399f8b1a608167c25/e2e/34d3d12d-039e-4496-a353-0d24175fbf15.md"; Disp = "e2e\34d3d12d-039e-4496-a353-0d24175fbf15.md" },
  @{ Cell = "B3"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06acceeb9e7bb8967abc64218bf0f6de1d20a0d1/e2e/546b8a45-a4fe-43f9-8570-96e9c4393b0d.md"; Disp = "e2e\546b8a45-a4fe-43f9-8570-96e9c4393b0d.md" },
  @{ Cell = "B4"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/274452995425a47c6d9ff5916fc705b1c4be371b/e2e/cfc2324b-6b69-48c4-8ec8-c64330098c47.md"; Disp = "e2e\cfc2324b-6b69-48c4-8ec8-c64330098c47.md" },
  @{ Cell = "B5"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5f8b17538436e140812b57399f8b1a608167c25/e2e/e0dcfb67-e9cf-4266-acbd-1203e67f0197.md"; Disp = "e2e\e0dcfb67-e9cf-4266-acbd-1203e67f0197.md" },
  @{ Cell = "B6"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37428f8f253c9e5417331a26628f7a5f243e298b/e2e/fb153bbd-eda5-4ccd-8490-bd45369ad1db.md"; Disp = "e2e\ca03050d-a2be-45e6-8869-2411d5b68e55.md" },
  @{ Cell = "B7"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/217ed6f1148f40541ee9baa8f73ee661c90a80aa/e2e/ca03050d-a2be-45e6-8869-2411d5b68e55.md"; Disp = "e2e\fb153bbd-eda5-4ccd-8490-bd45369ad1db.md" }
)
$ovLinks.Delete()
foreach ($t in $ovTargets) {
  $ovLinks.Add($ov.Range($t.Cell), $t.Url, "", "", $t.Disp) | Out-Null
}

# ---------------------------------------------------------------------------
# Helper data shared by the zh-cn and de-de sheets (columns A-P):
#   A Source File Name, B File Extension, C Status, D Source Path,
#   E Priority, F Content Duplicate, G Latest Handoff File,
#   H Latest Handoff Datetime, ...
# ---------------------------------------------------------------------------

# NOTE: this runtime's function calls only bind parameters positionally, so
# the helper below is invoked with plain positional arguments (no `-Name`
# syntax) to make sure the values actually reach the function body.
function Update-LangSheet($ws, $Row6Xlf, $Row6Datetime, $Row7Xlf, $Row7Datetime, $I4LinkUrl) {

  # Row 6 -> ca03050d, already "Ready for handoff".
  $ws.Range("A6").Value = "ca03050d-a2be-45e6-8869-2411d5b68e55.md"
  $ws.Range("C6").Value = "Ready for handoff"
  $ws.Range("G6").Value = $Row6Xlf
  $ws.Range("H6").Value = $Row6Datetime

  # Row 7 -> fb153bbd, newly "Ready for handoff" with a fresh handoff stamp.
  $ws.Range("A7").Value = "fb153bbd-eda5-4ccd-8490-bd45369ad1db.md"
  $ws.Range("C7").Value = "Ready for handoff"
  $ws.Range("G7").Value = $Row7Xlf
  $ws.Range("H7").Value = $Row7Datetime

  $links = $ws.Hyperlinks
  $targets = @(
    @{ Cell = "A2"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5f8b17538436e140812b57399f8b1a608167c25/e2e/34d3d12d-039e-4496-a353-0d24175fbf15.md"; Disp = "34d3d12d-039e-4496-a353-0d24175fbf15.md" },
    @{ Cell = "A3"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06acceeb9e7bb8967abc64218bf0f6de1d20a0d1/e2e/546b8a45-a4fe-43f9-8570-96e9c4393b0d.md"; Disp = "546b8a45-a4fe-43f9-8570-96e9c4393b0d.md" },
    @{ Cell = "A4"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/274452995425a47c6d9ff5916fc705b1c4be371b/e2e/cfc2324b-6b69-48c4-8ec8-c64330098c47.md"; Disp = "cfc2324b-6b69-48c4-8ec8-c64330098c47.md" },
    @{ Cell = "I4"; Url = $I4LinkUrl; Disp = "cfc2324b-6b69-48c4-8ec8-c64330098c47.md" },
    @{ Cell = "A5"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5f8b17538436e140812b57399f8b1a608167c25/e2e/e0dcfb67-e9cf-4266-acbd-1203e67f0197.md"; Disp = "e0dcfb67-e9cf-4266-acbd-1203e67f0197.md" },
    @{ Cell = "A6"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37428f8f253c9e5417331a26628f7a5f243e298b/e2e/fb153bbd-eda5-4ccd-8490-bd45369ad1db.md"; Disp = "ca03050d-a2be-45e6-8869-2411d5b68e55.md" },
    @{ Cell = "A7"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/217ed6f1148f40541ee9baa8f73ee661c90a80aa/e2e/ca03050d-a2be-45e6-8869-2411d5b68e55.md"; Disp = "fb153bbd-eda5-4ccd-8490-bd45369ad1db.md" }
  )
  $links.Delete()
  foreach ($t in $targets) {
    $links.Add($ws.Range($t.Cell), $t.Url, "", "", $t.Disp) | Out-Null
  }
}

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $zh `
  "ca03050d-a2be-45e6-8869-2411d5b68e55.3805c3a63486d4537a2cf0cc90ab88ec855b9f1e.zh-cn.xlf" `
  "2016-09-06 04:03:46" `
  "fb153bbd-eda5-4ccd-8490-bd45369ad1db.569916e2e055be12838f2459b3316bf9e1643a45.zh-cn.xlf" `
  "2016-09-06 04:10:56" `
  "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/26f7fe353530d0d3110f052c00493cd42a2e878a/e2e/cfc2324b-6b69-48c4-8ec8-c64330098c47.md"

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
Update-LangSheet $de `
  "ca03050d-a2be-45e6-8869-2411d5b68e55.3805c3a63486d4537a2cf0cc90ab88ec855b9f1e.de-de.xlf" `
  "2016-09-06 04:03:59" `
  "fb153bbd-eda5-4ccd-8490-bd45369ad1db.569916e2e055be12838f2459b3316bf9e1643a45.de-de.xlf" `
  "2016-09-06 04:11:11" `
  "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/66a055ddb29e2e23fa56372ca982c0e3bf09cae6/e2e/cfc2324b-6b69-48c4-8ec8-c64330098c47.md"

Write-Host "Done"
